$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '34.103.98'
$ws.Range("E2").Value = '  -0.13%  '

$ws.Range("D3").Value = '1.781.40'
$ws.Range("E3").Value = '  -0.61%  '

$ws.Range("E4").Value = '  +0.30%  '

Set-TextValue "D5" '225.56'
$ws.Range("E5").Value = '  -0.90%  '

$ws.Range("E6").Value = '  -0.21%  '

$ws.Range("E7").Value = '  +0.30%  '

$ws.Range("E8").Value = '  -1.26%  '

$ws.Range("E9").Value = '  -1.49%  '

$ws.Range("E10").Value = '  -0.33%  '

$ws.Range("E11").Value = '  +0.77%  '

$ws.Range("D12").Value = '2.038.87'
$ws.Range("E12").Value = '  -0.52%  '

$ws.Range("D13").Value = '1.790.45'
$ws.Range("E13").Value = '  -0.10%  '

Set-TextValue "D14" '10.91'
$ws.Range("E14").Value = '  -6.00%  '

$ws.Range("D16").Value = '34.094.83'

$ws.Range("E17").Value = '  -0.61%  '

Set-TextValue "D18" '67.51'
$ws.Range("E18").Value = '  -0.98%  '

Set-TextValue "D19" '245.24'
$ws.Range("E19").Value = '  +0.44%  '

$ws.Range("E20").Value = '  +1.15%  '

$ws.Range("E21").Value = '  +0.17%  '

Set-TextValue "D22" '10.88'
$ws.Range("E22").Value = '  +0.20%  '

$ws.Range("E23").Value = '  -0.39%  '

$ws.Range("E24").Value = '  -1.53%  '

Set-TextValue "D25" '161.93'
$ws.Range("E25").Value = '  -0.09%  '

$ws.Range("E26").Value = '  -0.93%  '

$ws.Range("E27").Value = '  -0.28%  '

$ws.Range("E28").Value = '  +0.32%  '

$ws.Range("E29").Value = '  +0.39%  '

$ws.Range("E30").Value = '  -1.41%  '

$ws.Range("E31").Value = '  -0.28%  '

$ws.Range("E32").Value = '  +0.87%  '

Set-TextValue "D33" '3.71'
$ws.Range("E33").Value = '  +1.91%  '

$ws.Range("E34").Value = '  -2.96%  '

$ws.Range("D35").Value = '1.446.24'

Set-TextValue "D36" '2.47'
$ws.Range("E36").Value = '  +5.40%  '

$ws.Range("E37").Value = '  -1.07%  '

Set-TextValue "D38" '0.0191'
$ws.Range("E38").Value = '  +0.49%  '

$ws.Range("E39").Value = '  -0.96%  '

$ws.Range("B40").Value = 'Aave'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue "D40" '81.14'
$ws.Range("E40").Value = '  +0.81%  '

$ws.Range("B41").Value = 'HuobiToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue "D41" '2.39'
$ws.Range("E41").Value = '  +1.36%  '

$ws.Range("E42").Value = '  +1.17%  '

Set-TextValue "D43" '0.912'
$ws.Range("E43").Value = '  -1.62%  '

$ws.Range("E44").Value = '  +1.83%  '

Set-TextValue "D45" '0.0520'
$ws.Range("E45").Value = '  +2.14%  '

$ws.Range("E46").Value = '  -0.30%  '

$ws.Range("E47").Value = '  +0.22%  '

$ws.Range("D48").Value = '1.938.61'
$ws.Range("E48").Value = '  -0.62%  '

$ws.Range("D49").Value = '0.0₆0133'
$ws.Range("E49").Value = '  -5.65%  '

Set-TextValue "D50" '104.67'
$ws.Range("E50").Value = '  -2.83%  '

$ws.Range("E51").Value = '  +0.30%  '
